$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.186.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.87%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.733.74'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.22%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '614.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '187.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.38%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.732.10'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.39%  '

$ws.Range("E8").Value = '  +1.63%  '

$ws.Range("E9").Value = '  -0.05%  '

$ws.Range("E10").Value = '  +0.73%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.162'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.37%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '56.70'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.97%  '

$ws.Range("E13").Value = '  -1.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.72'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.328.35'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.736.94'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.12'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.37'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.22%  '

$ws.Range("E19").Value = '  -0.24%  '

$ws.Range("E20").Value = '  +0.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.025.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '414.35'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.01%  '

$ws.Range("E23").Value = '  +2.57%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '89.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.32%  '

$ws.Range("E25").Value = '  -0.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.96%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.79%  '

$ws.Range("E29").Value = '  +0.30%  '

$ws.Range("E30").Value = '  +1.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.33'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.33'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -11.85%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.78'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.17%  '

$ws.Range("E34").Value = '  +3.06%  '

$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '66.53'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.91%  '

$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '44.61'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.99%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '614.52'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0863'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.55%  '

$ws.Range("E39").Value = '  +1.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.00%  '

$ws.Range("E41").Value = '  -0.26%  '

$ws.Range("E42").Value = '  +4.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.07'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.01%  '

$ws.Range("E44").Value = '  +0.33%  '

$ws.Range("E45").Value = '  +1.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.143'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.846.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.23'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.63%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.75%  '

$ws.Range("E50").Value = '  -15.48%  '

$ws.Range("E51").Value = '  -4.30%  '
